$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "1.8.2"
$wsMeta.Range("B8").Value = "2023-09-01T14:45:29-04:00"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("AJ1").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
